$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-08 Wednesday" "2025-01-14 Tuesday"

Replace-Text "78×51=" "39×85="
Replace-Text "59×11=" "92×74="
Replace-Text "22×48=" "97×47="
Replace-Text "77×62=" "78×21="
Replace-Text "78×38=" "84×30="
Replace-Text "37×48=" "95×97="
Replace-Text "85×50=" "57×19="
Replace-Text "79×71=" "71×55="
Replace-Text "16×86=" "62×71="
Replace-Text "81×12=" "99×54="
Replace-Text "69×73=" "14×50="
Replace-Text "96×98=" "41×47="
Replace-Text "64×95=" "86×76="
Replace-Text "89×44=" "82×11="
Replace-Text "36×83=" "56×63="
Replace-Text "47×42=" "45×56="
Replace-Text "29×36=" "23×13="
Replace-Text "21×31=" "99×88="
Replace-Text "22×65=" "93×54="
Replace-Text "22×88=" "31×33="
Replace-Text "12×87=" "11×64="
Replace-Text "31×85=" "35×93="
Replace-Text "73×56=" "56×34="
Replace-Text "46×49=" "45×87="
Replace-Text "44×89=" "51×48="
